$wb = $excel.ActiveWorkbook

# -- Status text: "Ready for handoff" -> "Handback transform failed" --
# Appears on Overview (row 8, zh-cn/de-de status columns) and on the
# per-language sheets' Status column (C8).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E8").Value = "Handback transform failed"
$wsOverview.Range("F8").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C8").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C8").Value = "Handback transform failed"

# -- New Error Detail text in the "Error Detail" column (P) for row 8 --
$wsZhCn.Range("P8").Value = "Handback file name: dcdc1txl.onq is different with handoff file name: 62c69bae-be57-4178-8ce6-43d74bae088d.18e6faf98997b83dff77e59d5b89d28f23b3ccf1.zh-cn."
$wsDeDe.Range("P8").Value = "Handback file name: dcdc1txl.onq is different with handoff file name: 62c69bae-be57-4178-8ce6-43d74bae088d.18e6faf98997b83dff77e59d5b89d28f23b3ccf1.de-de."

# -- Widen the "Error Detail" column (P, column 16) to fit the new text --
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
